$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from A10 (bold/bordered) to the newly created A17:A19 cells so they match existing column-A styling
$ws.Range("A10").Copy($ws.Range("A17"))
$ws.Range("A10").Copy($ws.Range("A18"))
$ws.Range("A10").Copy($ws.Range("A19"))

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.996862725871338
$ws.Range("D10").Value = 1.003936173464992
$ws.Range("E10").Value = 0.9946912582886986
$ws.Range("F10").Value = 0.996862725871338
$ws.Range("G10").Value = 0.9934932940071665
$ws.Range("H10").Value = 0.995918470698716
$ws.Range("I10").Value = 0.9935294117647059
$ws.Range("J10").Value = 1.003936173464992
$ws.Range("K10").Value = 0.9993137158768455
$ws.Range("L10").Value = 0.9980882208740917
$ws.Range("M10").Value = 0.9964052223492694

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9989556668137908
$ws.Range("D11").Value = 0.9800159008172991
$ws.Range("E11").Value = 0.9968120295885248
$ws.Range("F11").Value = 0.9989556668137908
$ws.Range("G11").Value = 0.9868734174947451
$ws.Range("H11").Value = 1.011543514075729
$ws.Range("I11").Value = 0.9985991305964876
$ws.Range("J11").Value = 0.9800159008172991
$ws.Range("K11").Value = 0.9884139652029119
$ws.Range("L11").Value = 0.9936848160083513
$ws.Range("M11").Value = 0.9954666098977626

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9988714597883372
$ws.Range("D12").Value = 0.9802644370539032
$ws.Range("E12").Value = 0.9967763015132033
$ws.Range("F12").Value = 0.9988714597883372
$ws.Range("G12").Value = 0.9870263211677281
$ws.Range("H12").Value = 1.011351989626171
$ws.Range("I12").Value = 0.9985411598059317
$ws.Range("J12").Value = 0.9802644370539032
$ws.Range("K12").Value = 0.9885203692835532
$ws.Range("L12").Value = 0.9936959145359452
$ws.Range("M12").Value = 0.9954719448258791

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9989365640913295
$ws.Range("D13").Value = 0.980074586621324
$ws.Range("E13").Value = 0.9967789019691861
$ws.Range("F13").Value = 0.9989365640913295
$ws.Range("G13").Value = 0.986906757659301
$ws.Range("H13").Value = 1.011494796296378
$ws.Range("I13").Value = 0.9985866750451996
$ws.Range("J13").Value = 0.980074586621324
$ws.Range("K13").Value = 0.988426744295255
$ws.Range("L13").Value = 0.9936816541932922
$ws.Range("M13").Value = 0.9954630469471196

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9887600000000007
$ws.Range("D14").Value = 1.016248
$ws.Range("E14").Value = 0.9975040000000008
$ws.Range("F14").Value = 0.9887600000000007
$ws.Range("G14").Value = 1.007488
$ws.Range("H14").Value = 0.9964759999999986
$ws.Range("I14").Value = 0.9914600000000005
$ws.Range("J14").Value = 1.016248
$ws.Range("K14").Value = 1.006876000000001
$ws.Range("L14").Value = 0.9978180000000008
$ws.Range("M14").Value = 0.9996560000000002

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.99
$ws.Range("D15").Value = 1.02
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.99
$ws.Range("G15").Value = 1.01
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 0.99
$ws.Range("J15").Value = 1.02
$ws.Range("K15").Value = 1.01
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.001666666666667

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9923900213248006
$ws.Range("D16").Value = 1.009666641919997
$ws.Range("E16").Value = 0.997863556300803
$ws.Range("F16").Value = 0.9923900213248006
$ws.Range("G16").Value = 1.004163989094396
$ws.Range("H16").Value = 0.998204544614402
$ws.Range("I16").Value = 0.9924549369856019
$ws.Range("J16").Value = 1.009666641919997
$ws.Range("K16").Value = 1.0037650991104
$ws.Range("L16").Value = 0.9980775602176003
$ws.Range("M16").Value = 0.9991239483733333

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9959560450014602
$ws.Range("D17").Value = 0.9966196613686858
$ws.Range("E17").Value = 0.9959556474424196
$ws.Range("F17").Value = 0.9959560450014602
$ws.Range("G17").Value = 0.9961045171712868
$ws.Range("H17").Value = 0.9963799365554631
$ws.Range("I17").Value = 0.9960756412583182
$ws.Range("J17").Value = 0.9966196613686858
$ws.Range("K17").Value = 0.9962876544055528
$ws.Range("L17").Value = 0.9961218497035065
$ws.Range("M17").Value = 0.9961819081329387

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.995433733333886
$ws.Range("D18").Value = 0.9944551982686333
$ws.Range("E18").Value = 0.9965967682354769
$ws.Range("F18").Value = 0.995433733333886
$ws.Range("G18").Value = 0.9958566478887693
$ws.Range("H18").Value = 0.9963585604965575
$ws.Range("I18").Value = 0.9971922848508404
$ws.Range("J18").Value = 0.9944551982686333
$ws.Range("K18").Value = 0.9955259832520551
$ws.Range("L18").Value = 0.9954798582929705
$ws.Range("M18").Value = 0.9959821988456938

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9962558235693746
$ws.Range("D19").Value = 0.9948720048294135
$ws.Range("E19").Value = 0.9963312775387843
$ws.Range("F19").Value = 0.9962558235693746
$ws.Range("G19").Value = 0.9953458371402598
$ws.Range("H19").Value = 0.9975014142485267
$ws.Range("I19").Value = 0.9962402831735806
$ws.Range("J19").Value = 0.9948720048294135
$ws.Range("K19").Value = 0.9956016411840989
$ws.Range("L19").Value = 0.9959287323767367
$ws.Range("M19").Value = 0.9960911067499899
